$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Vince",   "Hello",       "Hello",       "Vincent Martin",       "vince@nadin.one"),
    @("Chris95", "guest",       "guest",       "Christopher Perrins",  "chris@qa.com"),
    @("M4TT",    "chick3nWing", "chick3nWing", "matthewhunt",          "matt@qa.com"),
    @("Dev",     "D@T@",        "D@T@",        "Dev Gonsai",           "dev@qa.com"),
    @("matthewhunt", "chick3nWing", "chick3nWing", "Matttt",           "matthew.hunt@qa.com")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
